$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2,1).Value = "ECs"
$ws.Cells.Item(2,2).Value = "Col4a6"
$ws.Cells.Item(2,3).Value = "Cd93"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = 1
$ws.Cells.Item(2,6).Value = 0.3333333333333333
$ws.Cells.Item(2,7).Value = 0.012903
$ws.Cells.Item(2,8).Value = 0.038709
$ws.Cells.Item(2,9).Value = 0.007033395044444537
$ws.Cells.Item(2,10).Value = 0.007033395044444538
$ws.Cells.Item(2,11).Value = 3
$ws.Cells.Item(2,12).Value = 1
$ws.Cells.Item(2,13).Value = 83.91225566666667
$ws.Cells.Item(2,14).Value = 251.736767
$ws.Cells.Item(2,15).Value = 0.9556261553553385
$ws.Cells.Item(2,16).Value = 0.9556261553553385
$ws.Cells.Item(2,17).Value = 1.082719834867
$ws.Cells.Item(2,18).Value = 9.744478513803001
$ws.Cells.Item(2,19).Value = 0.006721296265417823
$ws.Cells.Item(2,20).Value = 0.006721296265417824

# Row 3
$ws.Cells.Item(3,1).Value = "ECs"
$ws.Cells.Item(3,2).Value = "Col4a6"
$ws.Cells.Item(3,3).Value = "Cd93"
$ws.Cells.Item(3,4).Value = "FAPs"
$ws.Cells.Item(3,5).Value = 1
$ws.Cells.Item(3,6).Value = 0.3333333333333333
$ws.Cells.Item(3,7).Value = 0.012903
$ws.Cells.Item(3,8).Value = 0.038709
$ws.Cells.Item(3,9).Value = 0.007033395044444537
$ws.Cells.Item(3,10).Value = 0.007033395044444538
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 0.3863573333333334
$ws.Cells.Item(3,14).Value = 1.159072
$ws.Cells.Item(3,15).Value = 0.00439999103960854
$ws.Cells.Item(3,16).Value = 0.00439999103960854
$ws.Cells.Item(3,17).Value = 0.004985168672000001
$ws.Cells.Item(3,18).Value = 0.04486651804800001
$ws.Cells.Item(3,19).Value = 0.00003094687517358307
$ws.Cells.Item(3,20).Value = 0.00003094687517358308

# Row 4
$ws.Cells.Item(4,1).Value = "ECs"
$ws.Cells.Item(4,2).Value = "Col4a6"
$ws.Cells.Item(4,3).Value = "Cd93"
$ws.Cells.Item(4,4).Value = "MuSCs"
$ws.Cells.Item(4,5).Value = 1
$ws.Cells.Item(4,6).Value = 0.3333333333333333
$ws.Cells.Item(4,7).Value = 0.012903
$ws.Cells.Item(4,8).Value = 0.038709
$ws.Cells.Item(4,9).Value = 0.007033395044444537
$ws.Cells.Item(4,10).Value = 0.007033395044444538
$ws.Cells.Item(4,11).Value = 3
$ws.Cells.Item(4,12).Value = 1
$ws.Cells.Item(4,13).Value = 3.510050666666667
$ws.Cells.Item(4,14).Value = 10.530152
$ws.Cells.Item(4,15).Value = 0.03997385360505296
$ws.Cells.Item(4,16).Value = 0.03997385360505297
$ws.Cells.Item(4,17).Value = 0.045290183752
$ws.Cells.Item(4,18).Value = 0.4076116537680001
$ws.Cells.Item(4,19).Value = 0.0002811519038531309
$ws.Cells.Item(4,20).Value = 0.000281151903853131

# Row 5
$ws.Cells.Item(5,1).Value = "FAPs"
$ws.Cells.Item(5,2).Value = "Col4a6"
$ws.Cells.Item(5,3).Value = "Cd93"
$ws.Cells.Item(5,4).Value = "ECs"
$ws.Cells.Item(5,5).Value = 3
$ws.Cells.Item(5,6).Value = 1
$ws.Cells.Item(5,7).Value = 1.07765
$ws.Cells.Item(5,8).Value = 3.23295
$ws.Cells.Item(5,9).Value = 0.5874244880760796
$ws.Cells.Item(5,10).Value = 0.5874244880760796
$ws.Cells.Item(5,11).Value = 3
$ws.Cells.Item(5,12).Value = 1
$ws.Cells.Item(5,13).Value = 83.91225566666667
$ws.Cells.Item(5,14).Value = 251.736767
$ws.Cells.Item(5,15).Value = 0.9556261553553385
$ws.Cells.Item(5,16).Value = 0.9556261553553385
$ws.Cells.Item(5,17).Value = 90.42804231918333
$ws.Cells.Item(5,18).Value = 813.85238087265
$ws.Cells.Item(5,19).Value = 0.5613582051017219
$ws.Cells.Item(5,20).Value = 0.5613582051017219

# Row 6
$ws.Cells.Item(6,1).Value = "FAPs"
$ws.Cells.Item(6,2).Value = "Col4a6"
$ws.Cells.Item(6,3).Value = "Cd93"
$ws.Cells.Item(6,4).Value = "FAPs"
$ws.Cells.Item(6,5).Value = 3
$ws.Cells.Item(6,6).Value = 1
$ws.Cells.Item(6,7).Value = 1.07765
$ws.Cells.Item(6,8).Value = 3.23295
$ws.Cells.Item(6,9).Value = 0.5874244880760796
$ws.Cells.Item(6,10).Value = 0.5874244880760796
$ws.Cells.Item(6,11).Value = 3
$ws.Cells.Item(6,12).Value = 1
$ws.Cells.Item(6,13).Value = 0.3863573333333334
$ws.Cells.Item(6,14).Value = 1.159072
$ws.Cells.Item(6,15).Value = 0.00439999103960854
$ws.Cells.Item(6,16).Value = 0.00439999103960854
$ws.Cells.Item(6,17).Value = 0.4163579802666667
$ws.Cells.Item(6,18).Value = 3.7472218224
$ws.Cells.Item(6,19).Value = 0.002584662483981384
$ws.Cells.Item(6,20).Value = 0.002584662483981384

# Row 7
$ws.Cells.Item(7,1).Value = "FAPs"
$ws.Cells.Item(7,2).Value = "Col4a6"
$ws.Cells.Item(7,3).Value = "Cd93"
$ws.Cells.Item(7,4).Value = "MuSCs"
$ws.Cells.Item(7,5).Value = 3
$ws.Cells.Item(7,6).Value = 1
$ws.Cells.Item(7,7).Value = 1.07765
$ws.Cells.Item(7,8).Value = 3.23295
$ws.Cells.Item(7,9).Value = 0.5874244880760796
$ws.Cells.Item(7,10).Value = 0.5874244880760796
$ws.Cells.Item(7,11).Value = 3
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 3.510050666666667
$ws.Cells.Item(7,14).Value = 10.530152
$ws.Cells.Item(7,15).Value = 0.03997385360505296
$ws.Cells.Item(7,16).Value = 0.03997385360505297
$ws.Cells.Item(7,17).Value = 3.782606100933334
$ws.Cells.Item(7,18).Value = 34.0434549084
$ws.Cells.Item(7,19).Value = 0.02348162049037639
$ws.Cells.Item(7,20).Value = 0.02348162049037639

# Row 8
$ws.Cells.Item(8,1).Value = "MuSCs"
$ws.Cells.Item(8,2).Value = "Col4a6"
$ws.Cells.Item(8,3).Value = "Cd93"
$ws.Cells.Item(8,4).Value = "ECs"
$ws.Cells.Item(8,5).Value = 3
$ws.Cells.Item(8,6).Value = 1
$ws.Cells.Item(8,7).Value = 0.7439806666666667
$ws.Cells.Item(8,8).Value = 2.231942
$ws.Cells.Item(8,9).Value = 0.4055421168794758
$ws.Cells.Item(8,10).Value = 0.4055421168794758
$ws.Cells.Item(8,11).Value = 3
$ws.Cells.Item(8,12).Value = 1
$ws.Cells.Item(8,13).Value = 83.91225566666667
$ws.Cells.Item(8,14).Value = 251.736767
$ws.Cells.Item(8,15).Value = 0.9556261553553385
$ws.Cells.Item(8,16).Value = 0.9556261553553385
$ws.Cells.Item(8,17).Value = 62.42909591239045
$ws.Cells.Item(8,18).Value = 561.8618632115141
$ws.Cells.Item(8,19).Value = 0.3875466539881988
$ws.Cells.Item(8,20).Value = 0.3875466539881988

# Row 9
$ws.Cells.Item(9,1).Value = "MuSCs"
$ws.Cells.Item(9,2).Value = "Col4a6"
$ws.Cells.Item(9,3).Value = "Cd93"
$ws.Cells.Item(9,4).Value = "FAPs"
$ws.Cells.Item(9,5).Value = 3
$ws.Cells.Item(9,6).Value = 1
$ws.Cells.Item(9,7).Value = 0.7439806666666667
$ws.Cells.Item(9,8).Value = 2.231942
$ws.Cells.Item(9,9).Value = 0.4055421168794758
$ws.Cells.Item(9,10).Value = 0.4055421168794758
$ws.Cells.Item(9,11).Value = 3
$ws.Cells.Item(9,12).Value = 1
$ws.Cells.Item(9,13).Value = 0.3863573333333334
$ws.Cells.Item(9,14).Value = 1.159072
$ws.Cells.Item(9,15).Value = 0.00439999103960854
$ws.Cells.Item(9,16).Value = 0.00439999103960854
$ws.Cells.Item(9,17).Value = 0.287442386424889
$ws.Cells.Item(9,18).Value = 2.586981477824
$ws.Cells.Item(9,19).Value = 0.001784381680453573
$ws.Cells.Item(9,20).Value = 0.001784381680453573

# Row 10
$ws.Cells.Item(10,1).Value = "MuSCs"
$ws.Cells.Item(10,2).Value = "Col4a6"
$ws.Cells.Item(10,3).Value = "Cd93"
$ws.Cells.Item(10,4).Value = "MuSCs"
$ws.Cells.Item(10,5).Value = 3
$ws.Cells.Item(10,6).Value = 1
$ws.Cells.Item(10,7).Value = 0.7439806666666667
$ws.Cells.Item(10,8).Value = 2.231942
$ws.Cells.Item(10,9).Value = 0.4055421168794758
$ws.Cells.Item(10,10).Value = 0.4055421168794758
$ws.Cells.Item(10,11).Value = 3
$ws.Cells.Item(10,12).Value = 1
$ws.Cells.Item(10,13).Value = 3.510050666666667
$ws.Cells.Item(10,14).Value = 10.530152
$ws.Cells.Item(10,15).Value = 0.03997385360505296
$ws.Cells.Item(10,16).Value = 0.03997385360505297
$ws.Cells.Item(10,17).Value = 2.611409835020445
$ws.Cells.Item(10,18).Value = 23.502688515184
$ws.Cells.Item(10,19).Value = 0.01621108121082344
$ws.Cells.Item(10,20).Value = 0.01621108121082345
